$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from B10 into the new B11 cell before setting its value,
# so that B11 matches the s="3" (wrap text, vertical top) style used by the rest of column B.
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A2").Value = @'
https://sebrae.com.br/sites/PortalSebrae/parceirodomunicipio/cidadeempreendedora/inovacao-e-sustentabilidade,88faf8298baaf810VgnVCM1000001b00320aRCRD
'@

$ws.Range("B2").Value = @'
Crie sua Conta Sebrae agora e tenha acesso gratuito a conteúdos e serviços imperdíveis. Já tenho uma conta.Entrar Fortalecendo o desenvolvimento local com soluções criativas e sustentáveis que geram riqueza, preservam o ambiente e melhoram a qualidade de vida  Quando se fala em inovação, muitas vezes associamos o termo às transformações tecnológicas e digitais, como inteligência artificial, robôs, carros autônomos e todas as novidades que vêm invadindo o nosso cotidiano. Mas o processo inovador não precisa da tecnologia para existir; ou melhor, é bem mais amplo do que qualquer ferramenta tecnológica. Inovação significa criar novas formas de fazer algo. É gerar riqueza e crescimento econômico por meio de ideias e soluções criativas.
'@

$ws.Range("A3").Value = @'
https://esginside.com.br/2025/09/30/tecnologias-limpas-e-transformacao-digital-impulsionam-solucoes-sustentaveis-no-brasil/
'@

$ws.Range("B3").Value = @'
Compartilhar O debate sobre o papel das tecnologias limpas na construção de um futuro sustentável reuniu três especialistas ligados ao IEEE: Teresa Cristina de Carvalho, professora da Escola Politécnica da USP; Otavio Chase, da Universidade Federal Rural da Amazônia (UFRA); e Vanessa Schramm, professora da Universidade Federal de Campina Grande (UFCG). Cada um trouxe perspectivas complementares sobre como ciência e inovação podem reduzir desigualdades, preservar recursos e ampliar o desenvolvimento econômico. Teresa Cristina de Carvalho apresentou projetos que unem preservação da floresta e tecnologia de ponta. Segundo ela, o programa Amazônia 4.0 busca instrumentalizar cadeias de produção sustentáveis com recursos da Indústria 4.0, promovendo rastreabilidade via IoT e QR Code e a distribuição equitativa dos benefícios para comunidades locais. Entre os exemplos, Teresa citou o projeto da cadeia do cacau e do cupuaçu, que utiliza estruturas geodésicas movidas a energia solar no Pará para capacitar indígenas e quilombolas na produção de chocolates de alto valor agregado no meio da flroesta.
'@

$ws.Range("A4").Value = @'
https://www.linkedin.com/pulse/o-papel-da-sustentabilidade-digital-como-inovar-sem-perder-qtisf/
'@

$ws.Range("B4").Value = @'
Agree & Join LinkedIn
      By clicking Continue to join or sign in, you agree to LinkedIn’s User Agreement, Privacy Policy, and Cookie Policy. Create your free account or sign in to continue your search
              or
      By clicking Continue to join or sign in, you agree to LinkedIn’s User Agreement, Privacy Policy, and Cookie Policy. New to LinkedIn? Join now
                          or
                    New to LinkedIn? Join now
      By clicking Continue to join or sign in, you agree to LinkedIn’s User Agreement, Privacy Policy, and Cookie Policy.
'@

$ws.Range("A5").Value = @'
https://www.gov.br/secom/pt-br/assuntos/noticias/2026/01/com-foco-em-inovacao-e-sustentabilidade-nova-industria-brasil-acelera-reindustrializacao-do-pais
'@

$ws.Range("B5").Value = @'
Notícias BALANÇO 2025 
    O programa Mover é um dos maiores símbolos desse novo momento, com a atração de investimentos e o fortalecimento da cadeia automotiva
    - Foto:
    Ricardo Stuckert/PR
   Ao longo de 2025, o Brasil consolidou a retomada dos investimentos na indústria, na inovação e na geração de empregos. A política da Nova Indústria Brasil (NIB) mostra, na prática, como o dinheiro público bem aplicado se transforma em crescimento econômico, mais oportunidades de trabalho e produtos melhores para a população. Os recursos do Plano Mais Produção, instrumento de financiamento da NIB, saltaram para R$ 643,3 bilhões no último ano. Desse total, já foram destinados R$ 588,4 bilhões, entre 2023 e 2025, para 406 mil projetos em todas as regiões, alinhados às seis missões da NIB. Dinheiro que fortaleceu fábricas, modernizou máquinas, estimulou novas tecnologias e ajudou empresas de todos os tamanhos a produzir mais, com eficiência e sustentabilidade.
'@

$ws.Range("A6").Value = @'
https://mittechreview.com.br/transformacao-digital-sustentavel-rumo-a-um-mundo-mais-responsavel-e-inovador/
'@

$ws.Range("B6").Value = @'
A Transformação Digital é uma oportunidade para empresas integrarem sustentabilidade em suas operações, promovendo um futuro responsável e competitivo. A combinação disso com tecnologias avançadas permite enfrentar desafios globais e prosperar em um mundo em evolução. A Transformação Digital não é apenas uma revolução tecnológica, mas também uma oportunidade crucial para as empresas alinharem suas operações e estratégias com a sustentabilidade, promovendo um futuro mais responsável e competitivo. A convergência desses dois conceitos não apenas impulsiona a eficiência, mas também abre portas para a inovação, a transparência e a responsabilidade, gerando diferenciais competitivos significativos. Sabemos que a TD permite que as empresas otimizem suas operações de maneira mais eficiente, podendo ser, ademais, sustentável.
'@

$ws.Range("A7").Value = @'
https://jornaldaindustria.com.br/tecnologia/rota-da-inovacao-aposta-em-solucoes-sustentaveis-e-digitais-para-transformar-o-brasil
'@

$ws.Range("B7").Value = @'
A jornada de inovação que percorrerá o Brasil em busca de tecnologias para transição ecológica e digital representa uma iniciativa estratégica que visa acelerar o desenvolvimento sustentável do país. A ação, coordenada pela Mobilização Empresarial pela Inovação (MEI), do Sistema Indústria, pretende mapear e conectar soluções inovadoras que possam impulsionar a transformação verde e digital em diversos setores produtivos. O projeto começa em julho e promete envolver ecossistemas de inovação regionais com foco em resultados práticos. A jornada de inovação que percorrerá o Brasil em busca de tecnologias para transição ecológica e digital passará por 14 estados, com eventos, visitas técnicas e encontros com startups, universidades, centros de pesquisa e empresas. A ideia é formar uma ponte entre os desafios industriais e as soluções tecnológicas capazes de oferecer mais eficiência, sustentabilidade e competitividade às empresas brasileiras.
'@

$ws.Range("A8").Value = @'
https://www.bv.com.br/bv-inspira/sustentabilidade/inovacao-e-sustentabilidade
'@

$ws.Range("B8").Value = @'
Descubra mais conteúdos  ESG Publicado em por Equipe BV Inspira Uma pesquisa feita pela Câmara Americana de Comércio (Amcham Brasil) – e divulgada pela CNN Brasil, mostrou que o percentual de companhias onde já são implementadas iniciativas de sustentabilidade passou de 71% para 76%. Esse é um avanço muito relevante e mostra como o tema vem ganhando espaço no mundo corporativo. O próximo passo é ampliar esse movimento, unindo inovação às práticas verdes para gerar cada vez mais valor de forma contínua. Mais do que nunca, existe a oportunidade de criar soluções que tragam retorno financeiro, social e ambiental de forma integrada, fortalecendo tanto os negócios quanto a sociedade. Neste artigo, você vai descobrir o que é inovação e sustentabilidade, por que elas se complementam, quais são os pilares da inovação sustentável e como colocar tudo isso em prática dentro da sua operação para uma gestão eficiente.
'@

$ws.Range("A9").Value = @'
https://blog.mdftechnology.com.br/ideias-10/
'@

$ws.Range("B9").Value = @'
Você já pensou em como as ideias inovadoras podem transformar nosso mundo? A tecnologia desempenha um papel vital na promoção da sustentabilidade, ajudando a criar soluções que beneficiam tanto o meio ambiente quanto a sociedade. Neste post, vamos discutir como as ideias mais brilhantes estão redefinindo o que significa viver de forma sustentável e como cada um de nós pode contribuir para essa mudança positiva. Sumário A tecnologia verde é uma área em crescimento que busca soluções que sejam ao mesmo tempo inovadoras e sustentáveis. Ao integrar princípios ecológicos em processos tecnológicos, podemos desenvolver soluções que limitam o impacto ambiental e promovem práticas de consumo responsáveis.
'@

$ws.Range("A10").Value = @'
https://blog.aevo.com.br/inovacao-em-sustentabilidade/
'@

$ws.Range("B10").Value = @'
A inovação em sustentabilidade abrange uma variedade de abordagens e práticas destinadas a criar soluções que atendam às necessidades do presente sem comprometer as gerações futuras. Ao integrar os princípios da sustentabilidade em todas as áreas da economia, a inovação nesse campo visa o equilíbrio entre os aspectos financeiros, sociais e ambientais do desenvolvimento. Nos últimos anos, a inovação em sustentabilidade assumiu um ritmo mais acelarado. Uma prova disso é o aumento considerável no número de patentes relacionadas aos Objetivos de Desenvolvimento Sustentável (ODS). Entre 2005 e 2020, os registros anuais saltaram de aproximadamente 1,8 milhão para mais de 5 milhões.
'@

$ws.Range("A11").Value = @'
https://centev.ufv.br/sustentabilidade-e-inovacao-caminhos-para-um-futuro-verde/
'@

$ws.Range("B11").Value = @'
Início » Sustentabilidade e Inovação: Caminhos para um Futuro Verde A sustentabilidade é um dos desafios mais urgentes do nosso tempo. Com a crescente preocupação com as mudanças climáticas, a escassez de recursos naturais e a degradação ambiental, a busca por soluções inovadoras se tornou essencial. Startups tecnológicas estão na vanguarda desta transformação, desenvolvendo produtos e serviços que não só atendem às necessidades do presente, mas também preservam o futuro. O papel das startups tecnológicas na promoção da sustentabilidade é crucial. Essas empresas estão utilizando a inovação para criar soluções que abordam problemas ambientais de maneiras novas e eficazes.
'@
